$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A14").Value = "2024-09-25T18:06:40Z"
$ws.Range("B14").Value = "temperature"
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "25"
$ws.Range("D14").Value = "N/A"
$ws.Range("E14").Value = "N/A"
$ws.Range("F14").Value = "N/A"
